$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.350.63'
$ws.Range('E2').Value = '  +1.91%  '
$ws.Range('D3').Value = '3.387.13'
$ws.Range('E3').Value = '  +1.67%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.70'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.74'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('E9').Value = '  +5.36%  '
$ws.Range('E10').Value = '  +1.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.46'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.68%  '
$ws.Range('E12').Value = '  +2.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '676.54'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.61'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.15%  '
$ws.Range('D15').Value = '3.931.40'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').Value = '69.422.42'
$ws.Range('E16').Value = '  +2.08%  '
$ws.Range('D17').Value = '3.437.61'
$ws.Range('E17').Value = '  +2.71%  '
$ws.Range('E18').Value = '  +1.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.61'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.27'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.902'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.44'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.15'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '103.51'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.71%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.66'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.06'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.69'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.12'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '556.31'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.59'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.14%  '
$ws.Range('E34').Value = '  +0.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.10'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.15%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = '3.686.75'
$ws.Range('E37').Value = '  -0.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.139'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +5.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.13'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.53%  '
$ws.Range('E40').Value = '  +2.53%  '
$ws.Range('E41').Value = '  +0.78%  '
$ws.Range('E42').Value = '  +3.31%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('E44').Value = '  +3.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.27'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.67'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.130'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.66%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.42'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +6.47%  '
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('E50').Value = '  +1.84%  '
$ws.Range('E51').Value = '  +2.79%  '
